$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.303.73"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "3.566.59"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "3.565.70"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +2.05%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "4.171.99"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "3.565.00"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("D17").Value = "66.330.47"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("E19").Value = "  +4.31%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.612"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").Value = "3.708.65"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").Value = "  +1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").Value = "3.560.99"
$ws.Range("E33").Value = "  +0.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  -5.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0851"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.890"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  +3.85%  "
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("E51").Value = "  +4.27%  "
